$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($rng.Find.Found) {
        $rng.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. Titre d'identité: Passeport -> Carte d'identité nationale
Replace-Text "Passeport" "Carte d'identité nationale"

# 2. Numéro de pièce : PP5456TRA -> AA-45467776-AQ
Replace-Text "PP5456TRA" "AA-45467776-AQ"

# 3. Date de délivrance : 02 mars 2020 -> 12 juillet 2023
Replace-Text "02 mars 2020" "12 juillet 2023"

# 4. Montant : 132 000 -> 430 000
Replace-Text "132 000" "430 000"

# 5. Date de signature : 04 décembre 2024 -> 11 décembre 2024
Replace-Text "04 décembre 2024" "11 décembre 2024"
